$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 340, shifting existing row 340 and below down by one.
$ws.Rows.Item(340).Insert()

# Populate the newly inserted row 340 with the new data record.
$ws.Cells.Item(340, 1).Value = 3
$ws.Cells.Item(340, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(340, 3).Value = "Coquimbo"
$ws.Cells.Item(340, 4).Value = 44855
$ws.Cells.Item(340, 5).Value = 5
$ws.Cells.Item(340, 6).Value = 100112009
$ws.Cells.Item(340, 7).Value = "Acelga"
$ws.Cells.Item(340, 8).Value = "Sin especificar"
$ws.Cells.Item(340, 9).Value = "Primera"
$ws.Cells.Item(340, 10).Value = 250
$ws.Cells.Item(340, 11).Value = 2000
$ws.Cells.Item(340, 12).Value = 2300
$ws.Cells.Item(340, 13).Value = 2144
$ws.Cells.Item(340, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(340, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(340, 16).Value = 357
$ws.Cells.Item(340, 17).Value = 6
$ws.Cells.Item(340, 18).Value = "Hortaliza"
